$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D, E and G keep their existing text (non-numeric) storage
# so round/percent-looking values are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "274.87"
$ws.Range("E2").Value = "-1.05%"
$ws.Range("G2").Value = "14"
$ws.Range("D3").Value = "27.41"
$ws.Range("E3").Value = "0.35%"
$ws.Range("G3").Value = "14"
$ws.Range("D4").Value = "4.788"
$ws.Range("E4").Value = "0.00%"
$ws.Range("G4").Value = "14"
$ws.Range("D5").Value = "0.06289"
$ws.Range("E5").Value = "-0.50%"
$ws.Range("G5").Value = "14"
$ws.Range("D6").Value = "6.920"
$ws.Range("E6").Value = "-0.27%"
$ws.Range("G6").Value = "14"
$ws.Range("D7").Value = "1.290"
$ws.Range("E7").Value = "36.13%"
$ws.Range("G7").Value = "14"
$ws.Range("D8").Value = "0.8719"
$ws.Range("E8").Value = "-1.06%"
$ws.Range("G8").Value = "14"
$ws.Range("D9").Value = "0.1511"
$ws.Range("E9").Value = "2.52%"
$ws.Range("G9").Value = "14"
$ws.Range("D10").Value = "0.05004"
$ws.Range("E10").Value = "-5.52%"
$ws.Range("G10").Value = "14"
$ws.Range("D11").Value = "0.07419"
$ws.Range("E11").Value = "1.28%"
$ws.Range("G11").Value = "14"
$ws.Range("D12").Value = "0.02892"
$ws.Range("E12").Value = "-7.56%"
$ws.Range("G12").Value = "14"
$ws.Range("D13").Value = "0.09059"
$ws.Range("E13").Value = "0.02%"
$ws.Range("G13").Value = "14"
$ws.Range("E14").Value = "0.30%"
$ws.Range("G14").Value = "14"
$ws.Range("D15").Value = "0.0006357"
$ws.Range("E15").Value = "1.44%"
$ws.Range("G15").Value = "14"
$ws.Range("D16").Value = "0.005870"
$ws.Range("E16").Value = "0.25%"
$ws.Range("G16").Value = "14"
$ws.Range("D17").Value = "3.451"
$ws.Range("E17").Value = "-0.40%"
$ws.Range("G17").Value = "14"
$ws.Range("D18").Value = "3.311"
$ws.Range("E18").Value = "-1.50%"
$ws.Range("G18").Value = "14"
$ws.Range("G19").Value = "14"
$ws.Range("E20").Value = "1.56%"
$ws.Range("G20").Value = "14"
$ws.Range("E21").Value = "-1.66%"
$ws.Range("G21").Value = "14"
$ws.Range("D22").Value = "3.898"
$ws.Range("E22").Value = "-0.11%"
$ws.Range("G22").Value = "14"
$ws.Range("D23").Value = "0.04369"
$ws.Range("E23").Value = "1.16%"
$ws.Range("G23").Value = "14"
$ws.Range("D24").Value = "0.001168"
$ws.Range("E24").Value = "-0.96%"
$ws.Range("G24").Value = "14"
$ws.Range("D25").Value = "0.003820"
$ws.Range("E25").Value = "6.41%"
$ws.Range("G25").Value = "14"
$ws.Range("D26").Value = "0.0001199"
$ws.Range("E26").Value = "-0.06%"
$ws.Range("G26").Value = "14"
$ws.Range("E27").Value = "-4.34%"
$ws.Range("G27").Value = "14"
$ws.Range("G28").Value = "14"
$ws.Range("G29").Value = "14"
$ws.Range("G30").Value = "14"
$ws.Range("G31").Value = "14"
$ws.Range("G32").Value = "14"
$ws.Range("G33").Value = "14"
$ws.Range("G34").Value = "14"
$ws.Range("G35").Value = "14"
$ws.Range("G36").Value = "14"
$ws.Range("G37").Value = "14"
$ws.Range("G38").Value = "14"
$ws.Range("G39").Value = "14"
$ws.Range("D40").Value = "0.04090"
$ws.Range("E40").Value = "1.86%"
$ws.Range("G40").Value = "14"
$ws.Range("D41").Value = "0.007039"
$ws.Range("E41").Value = "6.34%"
$ws.Range("G41").Value = "14"
$ws.Range("D42").Value = "0.1171"
$ws.Range("E42").Value = "0.89%"
$ws.Range("G42").Value = "14"
$ws.Range("D43").Value = "0.002019"
$ws.Range("E43").Value = "-13.73%"
$ws.Range("G43").Value = "14"
$ws.Range("D44").Value = "0.01121"
$ws.Range("E44").Value = "-9.69%"
$ws.Range("G44").Value = "14"
$ws.Range("D45").Value = "0.00005192"
$ws.Range("E45").Value = "-0.30%"
$ws.Range("G45").Value = "14"
$ws.Range("B46").Value = "BOLO"
$ws.Range("C46").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D46").Value = "1.486"
$ws.Range("E46").Value = "-37.51%"
$ws.Range("G46").Value = "14"
$ws.Range("B47").Value = "CoinbaseStockToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D47").Value = "0.02301"
$ws.Range("E47").Value = "2.23%"
$ws.Range("G47").Value = "14"
$ws.Range("G48").Value = "14"
$ws.Range("G49").Value = "14"
$ws.Range("G50").Value = "14"
$ws.Range("G51").Value = "14"
